$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.868.44'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '3.651.69'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '188.27'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.616'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.72%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.688'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.149'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '55.52'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000262'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -9.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('D14').Value = '4.238.89'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').Value = '3.655.22'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.64'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.37%  '
$ws.Range('D18').Value = '67.750.14'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.72%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '398.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.36'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '86.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.70'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.08'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.56'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '67.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.09'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.72%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '605.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '43.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.114'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.36%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.385'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.20%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.136'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0737'
$ws.Range('E41').Value = '  -16.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.84'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0417'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -11.60%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.134'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').Value = '2.736.70'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.75'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '142.01'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('E50').Value = '  -5.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -14.79%  '
